# Finalize version to deliver reports to user.
# Update CCX call-volume figures for rows 2-46 (columns B:J) to final values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(188,130,136,153,180,259,228,50,120)
    3 = @(41,34,36,28,46,61,57,21,30)
    4 = @(4,3,8,3,9,17,10,0,7)
    5 = @(10,4,4,1,4,4,13,5,9)
    6 = @(1,1,0,4,5,4,3,0,0)
    7 = @(8,11,7,10,14,13,12,11,3)
    8 = @(0,0,0,0,0,0,0,0,0)
    9 = @(9,9,6,4,5,7,4,3,7)
    10 = @(9,6,11,6,9,16,15,2,4)
    11 = @(51,31,24,40,41,59,56,4,25)
    12 = @(6,3,5,6,6,8,6,0,4)
    13 = @(4,2,0,4,6,5,1,4,0)
    14 = @(9,7,7,4,4,6,11,0,5)
    15 = @(8,9,3,8,4,11,12,0,2)
    16 = @(12,6,3,7,13,15,13,0,5)
    17 = @(11,3,6,10,8,13,13,0,9)
    18 = @(1,1,0,1,0,1,0,0,0)
    19 = @(14,14,16,10,10,13,23,8,14)
    20 = @(1,0,2,1,1,0,1,0,1)
    21 = @(0,0,0,0,0,0,0,0,0)
    22 = @(6,7,2,4,7,7,9,3,7)
    23 = @(2,5,7,4,1,6,10,5,6)
    24 = @(5,2,5,1,1,0,3,0,0)
    25 = @(38,6,27,28,39,54,44,3,5)
    26 = @(16,0,17,3,5,15,19,2,3)
    27 = @(7,3,1,10,15,10,9,1,2)
    28 = @(3,0,0,0,6,8,3,0,0)
    29 = @(3,2,3,15,13,14,4,0,0)
    30 = @(9,1,6,0,0,7,9,0,0)
    31 = @(0,0,0,0,0,0,0,0,0)
    32 = @(31,22,26,41,31,55,38,9,24)
    33 = @(3,3,3,9,4,8,5,0,1)
    34 = @(4,2,6,4,2,8,11,0,8)
    35 = @(3,5,3,3,6,6,2,1,2)
    36 = @(5,1,2,3,1,6,4,1,2)
    37 = @(0,0,1,4,1,5,0,1,5)
    38 = @(0,0,0,0,0,0,0,0,0)
    39 = @(10,7,8,15,15,19,10,4,3)
    40 = @(6,4,3,3,2,3,6,2,3)
    41 = @(13,23,7,6,13,17,10,5,22)
    42 = @(0,0,0,0,0,0,0,0,0)
    43 = @(6,5,1,1,8,4,2,0,0)
    44 = @(0,0,0,0,0,0,0,0,0)
    45 = @(4,3,4,3,5,10,5,2,4)
    46 = @(3,15,2,2,0,3,3,3,18)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}
